# Updated optimizer from bayesopt to random search - updated experiment
# and results of Lx parameter with new optimizer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths ---
# Note: Excel's ColumnWidth setter internally quantizes to whole-pixel
# (1/6 character-width) granularity before persisting to the sheet's <col>
# width attribute, exactly as real Excel does. The values below are the
# nearest representable widths to the target OOXML widths
# (12.7109375 / 14.7109375 / 11.7109375 "characters").
$ws.Columns.Item(3).ColumnWidth = 71/6
$ws.Columns.Item(5).ColumnWidth = 71/6
$ws.Columns.Item(6).ColumnWidth = 71/6
$ws.Columns.Item(7).ColumnWidth = 71/6
$ws.Columns.Item(8).ColumnWidth = 83/6
$ws.Columns.Item(9).ColumnWidth = 65/6

# --- Update data values ---

# Row 2
$ws.Range("B2").Value = 2.1421040753125338
$ws.Range("C2").Value = 0.22038107770705079
$ws.Range("D2").Value = 1.9107506265664163
$ws.Range("E2").Value = 0.9945833772911753
$ws.Range("F2").Value = 0.99728801120397281
$ws.Range("G2").Value = 1.0077798663324984
$ws.Range("H2").Value = 0.0054166227088247032
$ws.Range("I2").Value = 0.57649437598628916

# Row 3
$ws.Range("B3").Value = 1.1198306211201652
$ws.Range("C3").Value = 0.11520891163787707
$ws.Range("D3").Value = 1.0683914536974246
$ws.Range("E3").Value = 0.27180956736569606
$ws.Range("F3").Value = 0.52135359149592142
$ws.Range("G3").Value = 0.56349760216108913
$ws.Range("H3").Value = 0.72819043263430394
$ws.Range("I3").Value = 0.91154833997235962

# Row 4
$ws.Range("B4").Value = 4.2128375235700704
$ws.Range("C4").Value = 0.36825502828409712
$ws.Range("D4").Value = 3.3200000000000003
$ws.Range("E4").Value = 1.3231544574828158
$ws.Range("F4").Value = 1.1502845115373916
$ws.Range("G4").Value = 1.1723163841807913
$ws.Range("H4").Value = -0.32315445748281579
$ws.Range("I4").Value = 0.83980129657442193

$wb.Save()
